$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 204, pushing existing rows 204:292 down to 206:294.
$ws.Rows.Item(204).Insert()
$ws.Rows.Item(204).Insert()

# Populate the first new row (204) with a new "Primera" quality record.
$ws.Cells.Item(204, 1).Value  = 11
$ws.Cells.Item(204, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(204, 3).Value  = "Bíobío"
$ws.Cells.Item(204, 4).Value  = 44523
$ws.Cells.Item(204, 5).Value  = 8
$ws.Cells.Item(204, 6).Value  = 100112020
$ws.Cells.Item(204, 7).Value  = "Tomate"
$ws.Cells.Item(204, 8).Value  = "Larga vida"
$ws.Cells.Item(204, 9).Value  = "Primera"
$ws.Cells.Item(204, 10).Value = 800
$ws.Cells.Item(204, 11).Value = 12000
$ws.Cells.Item(204, 12).Value = 13000
$ws.Cells.Item(204, 13).Value = 12500
$ws.Cells.Item(204, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(204, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(204, 16).Value = 694
$ws.Cells.Item(204, 17).Value = 18
$ws.Cells.Item(204, 18).Value = "Hortaliza"

# Populate the second new row (205) with a new "Segunda" quality record.
$ws.Cells.Item(205, 1).Value  = 11
$ws.Cells.Item(205, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(205, 3).Value  = "Bíobío"
$ws.Cells.Item(205, 4).Value  = 44523
$ws.Cells.Item(205, 5).Value  = 8
$ws.Cells.Item(205, 6).Value  = 100112020
$ws.Cells.Item(205, 7).Value  = "Tomate"
$ws.Cells.Item(205, 8).Value  = "Larga vida"
$ws.Cells.Item(205, 9).Value  = "Segunda"
$ws.Cells.Item(205, 10).Value = 400
$ws.Cells.Item(205, 11).Value = 10000
$ws.Cells.Item(205, 12).Value = 10000
$ws.Cells.Item(205, 13).Value = 10000
$ws.Cells.Item(205, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(205, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(205, 16).Value = 556
$ws.Cells.Item(205, 17).Value = 18
$ws.Cells.Item(205, 18).Value = "Hortaliza"
